$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.691.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.369.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.79%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.46%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.378.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.410'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.970.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.897.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.385.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '364.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.525'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000122'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.40%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.51'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("E39").Value = '  -8.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.662.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0672'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '333.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0282'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.11%  '
